# Generate Report for Handoff
# Adds two new "Ready for handoff" entries (0aa7b477-... and 90675b55-...)
# ahead of the existing fb5b9956-... entry on all three sheets
# (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A (File Name) B (zh-cn) C (de-de) D (Latest Handoff Date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Drop every existing hyperlink on the sheet; inserting rows does not
# relocate hyperlink anchors in this engine, so they are rebuilt from
# scratch below once all of the row data is in its final place.
$wsOverview.Hyperlinks.Delete()

# Make room for the two new rows just above the fb5b9956 row (old row 5).
$wsOverview.Rows("5:6").Insert()

$wsOverview.Range("A5").Value2 = "0aa7b477-8e3e-4488-b6c6-334f45781283.md"
$wsOverview.Range("B5").Value2 = "Ready for handoff"
$wsOverview.Range("C5").Value2 = "Ready for handoff"
$wsOverview.Range("D5").Value2 = "2016-36-20 14:36:45"

$wsOverview.Range("A6").Value2 = "90675b55-f231-444b-877d-fec877996f6f.md"
$wsOverview.Range("B6").Value2 = "Ready for handoff"
$wsOverview.Range("C6").Value2 = "Ready for handoff"
$wsOverview.Range("D6").Value2 = "2016-36-20 14:36:45"

$wsOverview.Range("A7").Value2 = "fb5b9956-379f-4818-ba96-cff1dbe259bf.md"
$wsOverview.Range("B7").Value2 = "Ready for handoff"
$wsOverview.Range("C7").Value2 = "Ready for handoff"
$wsOverview.Range("D7").Value2 = "2016-34-20 14:34:48"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/00cf8624a6e621a6b687eff6d22cf84bfe575138/e2e/e5b8dcde-adf4-4289-95cd-6a9fca9a1ea9.md", "", "", "e5b8dcde-adf4-4289-95cd-6a9fca9a1ea9.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b7622e7f78fd19e79e5b510567af1e4237b928e/e2e/70d4f4c2-030b-4c72-bc2d-27fb913b736e.md", "", "", "70d4f4c2-030b-4c72-bc2d-27fb913b736e.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8b7622e7f78fd19e79e5b510567af1e4237b928e/e2e/b402c08e-499a-4bcf-9658-11a13547578c.md", "", "", "b402c08e-499a-4bcf-9658-11a13547578c.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5217553218def5d2ba4638a290de6dd36f515778/e2e/0aa7b477-8e3e-4488-b6c6-334f45781283.md", "", "", "0aa7b477-8e3e-4488-b6c6-334f45781283.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/28950466fcee07dd6d1f56175ab26b9feafa04a4/e2e/90675b55-f231-444b-877d-fec877996f6f.md", "", "", "90675b55-f231-444b-877d-fec877996f6f.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/72f4c1cd46b653c4f64405d84a9aee39bf9ac530/e2e/fb5b9956-379f-4818-ba96-cff1dbe259bf.md", "", "", "fb5b9956-379f-4818-ba96-cff1dbe259bf.md") | Out-Null

# ---------------------------------------------------------------------
# Locale detail sheets: zh-cn / de-de
# columns: A Source File Name, B File Extension, C Status,
#          D Latest Handoff File, E Latest Handoff Datetime,
#          H Handoff Reason, I Dependency From
# ---------------------------------------------------------------------
function Update-LocaleSheet($ws, $locale, $d5, $e5, $d6, $e6, $d7, $e7, $hoHash5, $hoHash6, $hoHash7) {

    $ws.Hyperlinks.Delete()
    $ws.Rows("5:6").Insert()

    $ws.Range("A5").Value2 = "0aa7b477-8e3e-4488-b6c6-334f45781283.md"
    $ws.Range("B5").Value2 = ".md"
    $ws.Range("C5").Value2 = "Ready for handoff"
    $ws.Range("D5").Value2 = $d5
    $ws.Range("E5").Value2 = $e5
    $ws.Range("H5").Value2 = "0001-01-01 00:00:00"
    $ws.Range("I5").Value2 = "Include"

    $ws.Range("A6").Value2 = "90675b55-f231-444b-877d-fec877996f6f.md"
    $ws.Range("B6").Value2 = ".md"
    $ws.Range("C6").Value2 = "Ready for handoff"
    $ws.Range("D6").Value2 = $d6
    $ws.Range("E6").Value2 = $e6
    $ws.Range("H6").Value2 = "0001-01-01 00:00:00"
    $ws.Range("I6").Value2 = "Include"

    $ws.Range("A7").Value2 = "fb5b9956-379f-4818-ba96-cff1dbe259bf.md"
    $ws.Range("B7").Value2 = ".md"
    $ws.Range("C7").Value2 = "Ready for handoff"
    $ws.Range("D7").Value2 = $d7
    $ws.Range("E7").Value2 = $e7
    $ws.Range("H7").Value2 = "0001-01-01 00:00:00"
    $ws.Range("I7").Value2 = "Include"

    # Row 2 (e5b8dcde)
    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/00cf8624a6e621a6b687eff6d22cf84bfe575138/e2e/e5b8dcde-adf4-4289-95cd-6a9fca9a1ea9.md", "", "", "e5b8dcde-adf4-4289-95cd-6a9fca9a1ea9.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/00cf8624a6e621a6b687eff6d22cf84bfe575138/e2e/e5b8dcde-adf4-4289-95cd-6a9fca9a1ea9.md", "", "", ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5b580010c60d01f7c5564082d635efb26f4445d8/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht/e5b8dcde-adf4-4289-95cd-6a9fca9a1ea9.1e8cf01ae400a2263b49353bc211542fe107a3d8.$locale.xlf", "", "", "e5b8dcde-adf4-4289-95cd-6a9fca9a1ea9.1e8cf01ae400a2263b49353bc211542fe107a3d8.$locale.xlf") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.$locale/blob/40c8d659556aae9d351bea4d00009b8b70c2de6f/e2e/e5b8dcde-adf4-4289-95cd-6a9fca9a1ea9.md", "", "", "e5b8dcde-adf4-4289-95cd-6a9fca9a1ea9.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1dc4891f34cc97dce0afb0027acf691f63eef909/ol-handback/OpenLocalizationTestOrg/oltest.$locale/ci/ht/e5b8dcde-adf4-4289-95cd-6a9fca9a1ea9.1e8cf01ae400a2263b49353bc211542fe107a3d8.$locale.xlf", "", "", "e5b8dcde-adf4-4289-95cd-6a9fca9a1ea9.1e8cf01ae400a2263b49353bc211542fe107a3d8.$locale.xlf") | Out-Null

    # Row 3 (70d4f4c2)
    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b7622e7f78fd19e79e5b510567af1e4237b928e/e2e/70d4f4c2-030b-4c72-bc2d-27fb913b736e.md", "", "", "70d4f4c2-030b-4c72-bc2d-27fb913b736e.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b7622e7f78fd19e79e5b510567af1e4237b928e/e2e/70d4f4c2-030b-4c72-bc2d-27fb913b736e.md", "", "", ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e69b58cb0d332aaa2f5b69a589c42937080bdec7/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht/70d4f4c2-030b-4c72-bc2d-27fb913b736e.9ca680c2afd801b401cb6c0610cf8a2c0dd701bf.$locale.xlf", "", "", "70d4f4c2-030b-4c72-bc2d-27fb913b736e.9ca680c2afd801b401cb6c0610cf8a2c0dd701bf.$locale.xlf") | Out-Null

    # Row 4 (b402c08e)
    $ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8b7622e7f78fd19e79e5b510567af1e4237b928e/e2e/b402c08e-499a-4bcf-9658-11a13547578c.md", "", "", "b402c08e-499a-4bcf-9658-11a13547578c.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/8b7622e7f78fd19e79e5b510567af1e4237b928e/e2e/b402c08e-499a-4bcf-9658-11a13547578c.md", "", "", ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e69b58cb0d332aaa2f5b69a589c42937080bdec7/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht/b402c08e-499a-4bcf-9658-11a13547578c.d3432c3e8d36f1b0d4d824fe43f34ce3ef057e24.$locale.xlf", "", "", "b402c08e-499a-4bcf-9658-11a13547578c.d3432c3e8d36f1b0d4d824fe43f34ce3ef057e24.$locale.xlf") | Out-Null

    # Row 5 (0aa7b477) - NEW
    $ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/$hoHash5/e2e/0aa7b477-8e3e-4488-b6c6-334f45781283.md", "", "", "0aa7b477-8e3e-4488-b6c6-334f45781283.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/$hoHash5/e2e/0aa7b477-8e3e-4488-b6c6-334f45781283.md", "", "", ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hoHash5/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht/0aa7b477-8e3e-4488-b6c6-334f45781283.5217553218def5d2ba4638a290de6dd36f515778.$locale.xlf", "", "", "0aa7b477-8e3e-4488-b6c6-334f45781283.5217553218def5d2ba4638a290de6dd36f515778.$locale.xlf") | Out-Null

    # Row 6 (90675b55) - NEW
    $ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$hoHash6/e2e/90675b55-f231-444b-877d-fec877996f6f.md", "", "", "90675b55-f231-444b-877d-fec877996f6f.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/$hoHash6/e2e/90675b55-f231-444b-877d-fec877996f6f.md", "", "", ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hoHash6/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht/90675b55-f231-444b-877d-fec877996f6f.28950466fcee07dd6d1f56175ab26b9feafa04a4.$locale.xlf", "", "", "90675b55-f231-444b-877d-fec877996f6f.28950466fcee07dd6d1f56175ab26b9feafa04a4.$locale.xlf") | Out-Null

    # Row 7 (fb5b9956) - moved down from row 5
    $ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/72f4c1cd46b653c4f64405d84a9aee39bf9ac530/e2e/fb5b9956-379f-4818-ba96-cff1dbe259bf.md", "", "", "fb5b9956-379f-4818-ba96-cff1dbe259bf.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/72f4c1cd46b653c4f64405d84a9aee39bf9ac530/e2e/fb5b9956-379f-4818-ba96-cff1dbe259bf.md", "", "", ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hoHash7/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht/fb5b9956-379f-4818-ba96-cff1dbe259bf.9098578b847812f099eddee2cfa549aae22e7add.$locale.xlf", "", "", "fb5b9956-379f-4818-ba96-cff1dbe259bf.9098578b847812f099eddee2cfa549aae22e7add.$locale.xlf") | Out-Null
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LocaleSheet $wsZhCn "zh-cn" `
    "0aa7b477-8e3e-4488-b6c6-334f45781283.5217553218def5d2ba4638a290de6dd36f515778.zh-cn.xlf" "2016-03-20 14:36:42" `
    "90675b55-f231-444b-877d-fec877996f6f.28950466fcee07dd6d1f56175ab26b9feafa04a4.zh-cn.xlf" "2016-03-20 14:36:42" `
    "fb5b9956-379f-4818-ba96-cff1dbe259bf.9098578b847812f099eddee2cfa549aae22e7add.zh-cn.xlf" "2016-03-20 14:34:45" `
    "5217553218def5d2ba4638a290de6dd36f515778" `
    "28950466fcee07dd6d1f56175ab26b9feafa04a4" `
    "8015b2593e758c2753384d6c443bfda91574c191"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LocaleSheet $wsDeDe "de-de" `
    "0aa7b477-8e3e-4488-b6c6-334f45781283.5217553218def5d2ba4638a290de6dd36f515778.de-de.xlf" "2016-03-20 14:36:45" `
    "90675b55-f231-444b-877d-fec877996f6f.28950466fcee07dd6d1f56175ab26b9feafa04a4.de-de.xlf" "2016-03-20 14:36:45" `
    "fb5b9956-379f-4818-ba96-cff1dbe259bf.9098578b847812f099eddee2cfa549aae22e7add.de-de.xlf" "2016-03-20 14:34:48" `
    "5217553218def5d2ba4638a290de6dd36f515778" `
    "28950466fcee07dd6d1f56175ab26b9feafa04a4" `
    "4c4ababb00b10edfcc4bbfa64fb46c998821968c"
